$d = $word.ActiveDocument

$p1 = $d.Paragraphs(1)

# Add a paragraph border (top/left/bottom/right, space=5pt from text) to paragraph 1.
$p1.Borders.DistanceFromTop = 5
$p1.Borders.DistanceFromLeft = 5
$p1.Borders.DistanceFromBottom = 5
$p1.Borders.DistanceFromRight = 5

# Widen the left indent of paragraph 1 from 6pt (120 twips) to 11.25pt (225 twips).
$p1.LeftIndent = 11.25

# Update the hidden ID placeholder text and drop the trailing space run.
$newId = "**ID__AFFARS_SUBPART_5318_000__ID**"
$oldLen = 31

$r1 = $d.Range(0, $oldLen)
$r1.Text = $newId

$spaceStart = $newId.Length
$r2 = $d.Range($spaceStart, $spaceStart + 1)
$r2.Delete()
